$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 288 - this shifts the existing rows
# 288..362 down to 289..363 (and grows the used range to A1:R363),
# mirroring Excel's native "insert row" behaviour (incl. carrying the
# date number-format down from the row above into the new D288 cell).
$ws.Rows(288).Insert()

# Populate the newly inserted row 288 with the new weekly record.
$ws.Range("A288").Value = 10
$ws.Range("B288").Value = "Vega Modelo de Temuco"
$ws.Range("C288").Value = "La Araucanía"
$ws.Range("D288").Value = 44855
$ws.Range("E288").Value = 9
$ws.Range("F288").Value = 100112001
$ws.Range("G288").Value = "Berenjena"
$ws.Range("H288").Value = "Sin especificar"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 50
$ws.Range("K288").Value = 15000
$ws.Range("L288").Value = 15000
$ws.Range("M288").Value = 15000
$ws.Range("N288").Value = "`$/caja 40 unidades"
$ws.Range("O288").Value = "Región Metropolitana"
$ws.Range("P288").Value = 375
$ws.Range("Q288").Value = 40
$ws.Range("R288").Value = "Hortaliza"
